# Update Distribution of Aberrant T Cells
# - Rename the Greek "Tgd" (T-gamma-delta) labels that used the Unicode
#   gamma/delta glyphs to the plain-ASCII "Tgd" spelling.
# - Because the cluster list is kept in alphabetical order, renaming
#   "Tγδ CD8+" / "Tγδ INSIG1+" to "Tgd CD8+" / "Tgd INSIG1+" moves those two
#   clusters from after "Trm IEL" to just before "Th17". The numeric data
#   for rows 28-32 is re-sorted to follow the renamed cluster labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple, standalone label renames (position in the sheet does not change).
$ws.Range("A2").Value = "Act. Tgd"
$ws.Range("A23").Value = "NK/Tgd"

# Capture the current (pre-edit) numeric data for the affected rows
# (28: Th17, 29: Tregs, 30: Trm IEL, 31: Tgd CD8+, 32: Tgd INSIG1+) before
# overwriting any of them, since the rows get reshuffled below.
# (.Value2 is used because it returns plain values/arrays for reads.)
$dataTh17   = $ws.Range("B28:V28").Value2
$dataTregs  = $ws.Range("B29:V29").Value2
$dataTrmIEL = $ws.Range("B30:V30").Value2
$dataTgdCD8 = $ws.Range("B31:V31").Value2
$dataTgdIns = $ws.Range("B32:V32").Value2

# New alphabetical order for rows 28-32:
#   28 Tgd CD8+
#   29 Tgd INSIG1+
#   30 Th17
#   31 Tregs
#   32 Trm IEL
$ws.Range("A28").Value = "Tgd CD8+"
$ws.Range("B28:V28").Value2 = $dataTgdCD8

$ws.Range("A29").Value = "Tgd INSIG1+"
$ws.Range("B29:V29").Value2 = $dataTgdIns

$ws.Range("A30").Value = "Th17"
$ws.Range("B30:V30").Value2 = $dataTh17

$ws.Range("A31").Value = "Tregs"
$ws.Range("B31:V31").Value2 = $dataTregs

$ws.Range("A32").Value = "Trm IEL"
$ws.Range("B32:V32").Value2 = $dataTrmIEL
